$wb = $excel.ActiveWorkbook

# ---- "axis" sheet: selection cursor moved C4 -> A9 ----
$wsAxis = $wb.Worksheets.Item("axis")
[void]$wsAxis.Activate()
[void]$wsAxis.Range("A9").Select()

# ---- "structure" sheet: term-tracking updates ----
$ws = $wb.Worksheets.Item("structure")
[void]$ws.Activate()

# Insert a new leading column to hold a per-term "status" note.
$ws.Columns("A:A").Insert()

# Add the new "talus head" term as row 9 (entered in this order so the
# shared-string table fills in the same sequence as the source edit).
$ws.Range("B9").Value = "talus head"
$ws.Range("A9").Value = "made synnonym request on phenotype-ext"
$ws.Range("C9").Value = "astragalus head"

# Match row 9's label cell formatting to its sibling label cells.
$ws.Range("B2").Copy()
[void]$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new "status" column: header + per-row progress notes.
$ws.Range("A1").Value = "status"
$ws.Range("A2").Value = "on fovt-talus; made pull request"
$ws.Range("A3").Value = "on fovt-talus; made pull request"
$ws.Range("A4").Value = "on fovt-talus; made pull request"
$ws.Range("A5").Value = "on fovt-talus; made pull request"
$ws.Range("A6").Value = "on fovt-talus; made pull request"
$ws.Range("A7").Value = "on fovt-talus; made pull request"
$ws.Range("A8").Value = "on fovt-talus; made pull request"
$ws.Range("A10").Value = "on fovt-talus; made pull request"

# Leave the selection/view where the author left it.
[void]$ws.Range("B18").Select()
